$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
#
# The underlying data source re-sorted some fixtures; for four pairs of
# adjacent rows the match data (columns B..AC) needs to be swapped between
# the two rows while each row keeps its own running index in column A.

function Swap-Rows([int]$row1, [int]$row2) {
    $r1 = $ws.Range("B$row1`:AC$row1")
    $r2 = $ws.Range("B$row2`:AC$row2")
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

Swap-Rows 139 140
Swap-Rows 177 178
Swap-Rows 195 196
Swap-Rows 231 232

# A handful of odds were independently corrected on other fixtures.
$ws.Range("N244").Value = 2.05
$ws.Range("P244").Value = 3.75
$ws.Range("Q244").Value = -0.5
$ws.Range("R244").Value = 2.08
$ws.Range("S244").Value = 1.82

$ws.Range("O245").Value = 3.2
$ws.Range("P245").Value = 3.3
$ws.Range("U245").Value = 2.025
$ws.Range("V245").Value = 1.825

$ws.Range("R246").Value = 1.82
$ws.Range("S246").Value = 2.08

$ws.Range("R247").Value = 2.01
$ws.Range("S247").Value = 1.89
$ws.Range("U247").Value = 1.85
$ws.Range("V247").Value = 2

$ws.Range("U248").Value = 1.85
$ws.Range("V248").Value = 2

$ws.Range("R249").Value = 1.89
$ws.Range("S249").Value = 2.01

$ws.Range("R251").Value = 1.88
$ws.Range("S251").Value = 2.02
$ws.Range("U251").Value = 1.925
$ws.Range("V251").Value = 1.925

$ws.Range("R252").Value = 2.11
$ws.Range("S252").Value = 1.79

$wb.Save()
